$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) and a single data row
# (row 2) for "Mujeeb Ur Rahman". The commit duplicates that data row,
# appending it as a new row 3 (same venue/date/result/teams/stats).
# Copy row 2 verbatim into row 3 so the values/types (including the
# numbers that are stored as text) are preserved exactly.
$ws.Range("A2:K2").Copy()
$ws.Range("A3").PasteSpecial()
